# Generate Report for Archive
# The localization status for the two handed-back files moved on from
# "Ready for handoff" to "In Translation". Update every sheet that surfaces
# that status (the Overview roll-up plus each per-locale detail sheet), then
# re-fit the status column(s) now that the new text is shorter.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: one status column per locale (E = zh-cn, F = de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- Per-locale detail sheets: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Re-fit the status columns now that "In Translation" is shorter than
#     "Ready for handoff" ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
